# "added slide titles for all slides"
#
# The deck starts with just the title slide ("make"). This adds the
# eight follow-on slides for the session, each using the "Title and
# Content" layout (the same layout PowerPoint picks when you type a new
# title in Outline view), and sets each slide's title text. The slides
# are inserted directly in their final, presented order.

$p = $ppt.ActivePresentation

$titles = @(
    "More about compiling",
    "Build systems",
    "Make",
    "Pattern rules",
    "Special variables",
    "Dependencies",
    "Targets",
    "Exercise: re-do bash workflow in make"
)

$idx = 2
foreach ($title in $titles) {
    $slide = $p.Slides.Add($idx, 2)
    $titleShape = $slide.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Text = $title

    # The longest title overflows the placeholder at full size; PowerPoint
    # shrinks it to fit automatically.
    if ($title.Length -gt 35) {
        $titleShape.TextFrame.AutoSize = 2
    }

    $idx = $idx + 1
}
